$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.456.51'
$ws.Range("E2").Value = '  +1.81%  '

# Row 3
$ws.Range("D3").Value = '2.612.14'
$ws.Range("E3").Value = '  +1.71%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.79'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.35%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.20'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.71%  '

# Row 7
$ws.Range("E7").Value = '  +0.15%  '

# Row 8
$ws.Range("E8").Value = '  +1.41%  '

# Row 9
$ws.Range("D9").Value = '2.621.75'
$ws.Range("E9").Value = '  +1.56%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.48'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.31%  '

# Row 11
$ws.Range("E11").Value = '  +4.30%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.338'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.97%  '

# Row 13
$ws.Range("E13").Value = '  +1.92%  '

# Row 14
$ws.Range("D14").Value = '3.067.68'
$ws.Range("E14").Value = '  +1.66%  '

# Row 15
$ws.Range("D15").Value = '59.331.76'
$ws.Range("E15").Value = '  +1.74%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.59'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.21%  '

# Row 17
$ws.Range("D17").Value = '2.615.06'
$ws.Range("E17").Value = '  +1.89%  '

# Row 18
$ws.Range("E18").Value = '  +2.71%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '347.06'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.77%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.36'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.22%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.16'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.19%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.41'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.76%  '

# Row 23
$ws.Range("E23").Value = '  -0.23%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.11'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.04%  '

# Row 25
$ws.Range("E25").Value = '  +1.28%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.410'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.58%  '

# Row 27
$ws.Range("E27").Value = '  +0.27%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.22'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.78%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0749'
$ws.Range("E29").Value = '  +7.92%  '

# Row 30
$ws.Range("E30").Value = '  +0.13%  '

# Row 31
$ws.Range("E31").Value = '  +6.25%  '

# Row 32
$ws.Range("E32").Value = '  +0.43%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.90'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.77%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.03'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.54%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.01'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.38%  '

# Row 36
$ws.Range("E36").Value = '  +1.40%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.97'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.843'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.50%  '

# Row 39
$ws.Range("E39").Value = '  +3.55%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.842'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.25%  '

# Row 41
$ws.Range("E41").Value = '  +2.54%  '

# Row 42
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '277.27'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.73%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.23%  '

# Row 44
$ws.Range("E44").Value = '  +2.95%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.76'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.31%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0963'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.36%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0525'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.43%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.61'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.84%  '

# Row 49
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.950.15'
$ws.Range("E49").Value = '  -0.49%  '

# Row 50
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0224'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.16%  '

# Row 51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.40'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.69%  '
